$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.444.83'
$ws.Range('E2').Value = '  -1.21%  '
$ws.Range('D3').Value = '2.312.81'
$ws.Range('E3').Value = '  -1.58%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '511.04'
$ws.Range('E5').Value = '  -1.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.07'
$ws.Range('E6').Value = '  -2.88%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  -1.18%  '
$ws.Range('E9').Value = '  -3.62%  '
$ws.Range('E10').Value = '  -0.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.23'
$ws.Range('E11').Value = '  -0.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.335'
$ws.Range('E12').Value = '  -2.28%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '2.727.06'
$ws.Range('E13').Value = '  -1.08%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.39'
$ws.Range('E14').Value = '  -2.20%  '
$ws.Range('D15').Value = '56.405.08'
$ws.Range('E15').Value = '  -1.01%  '
$ws.Range('E16').Value = '  -2.11%  '
$ws.Range('D17').Value = '2.323.66'
$ws.Range('E17').Value = '  -0.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.33'
$ws.Range('E18').Value = '  -1.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '326.19'
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('E20').Value = '  -2.52%  '
$ws.Range('E21').Value = '  +1.46%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.22'
$ws.Range('E23').Value = '  +0.63%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.54'
$ws.Range('E24').Value = '  +7.53%  '
$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.163'
$ws.Range('E25').Value = '  -1.11%  '
$ws.Range('E26').Value = '  +2.54%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '167.26'
$ws.Range('E28').Value = '  -1.09%  '
$ws.Range('E29').Value = '  -3.06%  '
$ws.Range('D30').Value = '0.0₃0714'
$ws.Range('E30').Value = '  -4.31%  '
$ws.Range('E31').Value = '  -2.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.25'
$ws.Range('E32').Value = '  -0.39%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.24'
$ws.Range('E35').Value = '  -0.63%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.90'
$ws.Range('E36').Value = '  -3.57%  '
$ws.Range('E37').Value = '  -4.99%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '38.54'
$ws.Range('E38').Value = '  +1.63%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.55'
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '149.19'
$ws.Range('E40').Value = '  +7.66%  '
$ws.Range('E41').Value = '  -1.89%  '
$ws.Range('E42').Value = '  -1.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '275.24'
$ws.Range('E43').Value = '  -0.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.97'
$ws.Range('E44').Value = '  -5.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0925'
$ws.Range('E45').Value = '  -0.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0493'
$ws.Range('E46').Value = '  -2.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.552'
$ws.Range('E47').Value = '  -2.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.06'
$ws.Range('E48').Value = '  +0.96%  '
$ws.Range('E49').Value = '  -1.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.93'
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '10.98'
$ws.Range('E51').Value = '  +0.65%  '
